$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Variable"): update row 5 in place, then clone its formatting into rows 6-20 ---
$ws1.Range("A5:K5").Copy($ws1.Range("A6:K20"))

$ws1.Range("A5").Value = "CREATE/MODIFY"
$ws1.Range("B5").Value = "BE_IND_1"
$ws1.Range("C5").Value = "BE_IND_1"
$ws1.Range("D5").Value = "BE_IND_1"
$ws1.Range("E5").Value = "IND_1"
$ws1.Range("F5").Value = "1 - Days past due"
$ws1.Range("G5").Value = 1
$ws1.Range("H5").Value = $false
$ws1.Range("I5").Value = $false
$ws1.Range("J5").Value = "CONTINUOUS"
$ws1.Range("K5").Value = "true"

$ws1.Range("A6").Value = "CREATE/MODIFY"
$ws1.Range("B6").Value = "BE_IND_2"
$ws1.Range("C6").Value = "BE_IND_2"
$ws1.Range("D6").Value = "BE_IND_2"
$ws1.Range("E6").Value = "IND_2"
$ws1.Range("F6").Value = "2 - Past Due > 90"
$ws1.Range("G6").Value = 2
$ws1.Range("H6").Value = $false
$ws1.Range("I6").Value = $false
$ws1.Range("J6").Value = "CATEGORICAL"
$ws1.Range("K6").Value = "true"

$ws1.Range("A7").Value = "CREATE/MODIFY"
$ws1.Range("B7").Value = "BE_IND_3"
$ws1.Range("C7").Value = "BE_IND_3"
$ws1.Range("D7").Value = "BE_IND_3"
$ws1.Range("E7").Value = "IND_3"
$ws1.Range("F7").Value = "3 - Blocked accounts"
$ws1.Range("G7").Value = 3
$ws1.Range("H7").Value = $false
$ws1.Range("I7").Value = $false
$ws1.Range("J7").Value = "CATEGORICAL"
$ws1.Range("K7").Value = "true"

$ws1.Range("A8").Value = "CREATE/MODIFY"
$ws1.Range("B8").Value = "BE_IND_7"
$ws1.Range("C8").Value = "BE_IND_7"
$ws1.Range("D8").Value = "BE_IND_7"
$ws1.Range("E8").Value = "IND_7"
$ws1.Range("F8").Value = "7 - Overdue amount/exposure amount"
$ws1.Range("G8").Value = 7
$ws1.Range("H8").Value = $false
$ws1.Range("I8").Value = $false
$ws1.Range("J8").Value = "CONTINUOUS"
$ws1.Range("K8").Value = "true"

$ws1.Range("A9").Value = "CREATE/MODIFY"
$ws1.Range("B9").Value = "BE_IND_8"
$ws1.Range("C9").Value = "BE_IND_8"
$ws1.Range("D9").Value = "BE_IND_8"
$ws1.Range("E9").Value = "IND_8"
$ws1.Range("F9").Value = "8 - Account turnover oscillation"
$ws1.Range("G9").Value = 8
$ws1.Range("H9").Value = $false
$ws1.Range("I9").Value = $false
$ws1.Range("J9").Value = "CONTINUOUS"
$ws1.Range("K9").Value = "true"

$ws1.Range("A10").Value = "CREATE/MODIFY"
$ws1.Range("B10").Value = "BE_IND_9"
$ws1.Range("C10").Value = "BE_IND_9"
$ws1.Range("D10").Value = "BE_IND_9"
$ws1.Range("E10").Value = "IND_9"
$ws1.Range("F10").Value = "9 - Delta turnover"
$ws1.Range("G10").Value = 9
$ws1.Range("H10").Value = $false
$ws1.Range("I10").Value = $false
$ws1.Range("J10").Value = "CATEGORICAL"
$ws1.Range("K10").Value = "true"

$ws1.Range("A11").Value = "CREATE/MODIFY"
$ws1.Range("B11").Value = "BE_IND_14"
$ws1.Range("C11").Value = "BE_IND_14"
$ws1.Range("D11").Value = "BE_IND_14"
$ws1.Range("E11").Value = "IND_14"
$ws1.Range("F11").Value = "14 - Number of business current accounts "
$ws1.Range("G11").Value = 14
$ws1.Range("H11").Value = $false
$ws1.Range("I11").Value = $false
$ws1.Range("J11").Value = "CONTINUOUS"
$ws1.Range("K11").Value = "true"

$ws1.Range("A12").Value = "CREATE/MODIFY"
$ws1.Range("B12").Value = "BE_IND_16"
$ws1.Range("C12").Value = "BE_IND_16"
$ws1.Range("D12").Value = "BE_IND_16"
$ws1.Range("E12").Value = "IND_16"
$ws1.Range("F12").Value = "16 - Default"
$ws1.Range("G12").Value = 16
$ws1.Range("H12").Value = $false
$ws1.Range("I12").Value = $false
$ws1.Range("J12").Value = "CATEGORICAL"
$ws1.Range("K12").Value = "true"

$ws1.Range("A13").Value = "CREATE/MODIFY"
$ws1.Range("B13").Value = "BE_IND_34"
$ws1.Range("C13").Value = "BE_IND_34"
$ws1.Range("D13").Value = "BE_IND_34"
$ws1.Range("E13").Value = "IND_34"
$ws1.Range("F13").Value = "34 - Negative own equity"
$ws1.Range("G13").Value = 34
$ws1.Range("H13").Value = $false
$ws1.Range("I13").Value = $false
$ws1.Range("J13").Value = "CATEGORICAL"
$ws1.Range("K13").Value = "true"

$ws1.Range("A14").Value = "CREATE/MODIFY"
$ws1.Range("B14").Value = "BE_IND_35"
$ws1.Range("C14").Value = "BE_IND_35"
$ws1.Range("D14").Value = "BE_IND_35"
$ws1.Range("E14").Value = "IND_35"
$ws1.Range("F14").Value = "35 - Delta equity"
$ws1.Range("G14").Value = 35
$ws1.Range("H14").Value = $false
$ws1.Range("I14").Value = $false
$ws1.Range("J14").Value = "CONTINUOUS"
$ws1.Range("K14").Value = "true"

$ws1.Range("A15").Value = "CREATE/MODIFY"
$ws1.Range("B15").Value = "BE_IND_40"
$ws1.Range("C15").Value = "BE_IND_40"
$ws1.Range("D15").Value = "BE_IND_40"
$ws1.Range("E15").Value = "IND_40"
$ws1.Range("F15").Value = "40 - Loan to value ratio"
$ws1.Range("G15").Value = 40
$ws1.Range("H15").Value = $false
$ws1.Range("I15").Value = $false
$ws1.Range("J15").Value = "CONTINUOUS"
$ws1.Range("K15").Value = "true"

$ws1.Range("A16").Value = "CREATE/MODIFY"
$ws1.Range("B16").Value = "BE_IND_44"
$ws1.Range("C16").Value = "BE_IND_44"
$ws1.Range("D16").Value = "BE_IND_44"
$ws1.Range("E16").Value = "IND_44"
$ws1.Range("F16").Value = "44 - Past due amount"
$ws1.Range("G16").Value = 44
$ws1.Range("H16").Value = $false
$ws1.Range("I16").Value = $false
$ws1.Range("J16").Value = "CONTINUOUS"
$ws1.Range("K16").Value = "true"

$ws1.Range("A17").Value = "CREATE/MODIFY"
$ws1.Range("B17").Value = "BE_IND_48"
$ws1.Range("C17").Value = "BE_IND_48"
$ws1.Range("D17").Value = "BE_IND_48"
$ws1.Range("E17").Value = "IND_48"
$ws1.Range("F17").Value = "48 - Debt Service Coverage Ratio"
$ws1.Range("G17").Value = 48
$ws1.Range("H17").Value = $false
$ws1.Range("I17").Value = $false
$ws1.Range("J17").Value = "CONTINUOUS"
$ws1.Range("K17").Value = "true"

$ws1.Range("A18").Value = "CREATE/MODIFY"
$ws1.Range("B18").Value = "BE_IND_51"
$ws1.Range("C18").Value = "BE_IND_51"
$ws1.Range("D18").Value = "BE_IND_51"
$ws1.Range("E18").Value = "IND_51"
$ws1.Range("F18").Value = "51 - Overdraft"
$ws1.Range("G18").Value = 51
$ws1.Range("H18").Value = $false
$ws1.Range("I18").Value = $false
$ws1.Range("J18").Value = "CATEGORICAL"
$ws1.Range("K18").Value = "true"

$ws1.Range("A19").Value = "CREATE/MODIFY"
$ws1.Range("B19").Value = "BE_IND_55"
$ws1.Range("C19").Value = "BE_IND_55"
$ws1.Range("D19").Value = "BE_IND_55"
$ws1.Range("E19").Value = "IND_55"
$ws1.Range("F19").Value = "55 - Forborne NPE"
$ws1.Range("G19").Value = 55
$ws1.Range("H19").Value = $false
$ws1.Range("I19").Value = $false
$ws1.Range("J19").Value = "CONTINUOUS"
$ws1.Range("K19").Value = "true"

$ws1.Range("A20").Value = "CREATE/MODIFY"
$ws1.Range("B20").Value = "BE_IND_56"
$ws1.Range("C20").Value = "BE_IND_56"
$ws1.Range("D20").Value = "BE_IND_56"
$ws1.Range("E20").Value = "IND_56"
$ws1.Range("F20").Value = "56 - Outstanding + overdue/Approved amount for loans"
$ws1.Range("G20").Value = 56
$ws1.Range("H20").Value = $false
$ws1.Range("I20").Value = $false
$ws1.Range("J20").Value = "CONTINUOUS"
$ws1.Range("K20").Value = "true"

# --- Sheet2 ("r Variable_DataType"): update row 5 in place, then clone its formatting into rows 6-20 ---
$ws2.Range("A5:F5").Copy($ws2.Range("A6:F20"))

$ws2.Range("A5").Value = "CREATE/MODIFY"
$ws2.Range("B5").Value = "BE_IND_1_REAL"
$ws2.Range("C5").Value = "BE_IND_1_REAL"
$ws2.Range("E5").Value = "BE_IND_1"
$ws2.Range("F5").Value = "REAL"

$ws2.Range("A6").Value = "CREATE/MODIFY"
$ws2.Range("B6").Value = "BE_IND_2_STRING"
$ws2.Range("C6").Value = "BE_IND_2_STRING"
$ws2.Range("E6").Value = "BE_IND_2"
$ws2.Range("F6").Value = "STRING"

$ws2.Range("A7").Value = "CREATE/MODIFY"
$ws2.Range("B7").Value = "BE_IND_3_STRING"
$ws2.Range("C7").Value = "BE_IND_3_STRING"
$ws2.Range("E7").Value = "BE_IND_3"
$ws2.Range("F7").Value = "STRING"

$ws2.Range("A8").Value = "CREATE/MODIFY"
$ws2.Range("B8").Value = "BE_IND_7_REAL"
$ws2.Range("C8").Value = "BE_IND_7_REAL"
$ws2.Range("E8").Value = "BE_IND_7"
$ws2.Range("F8").Value = "REAL"

$ws2.Range("A9").Value = "CREATE/MODIFY"
$ws2.Range("B9").Value = "BE_IND_8_REAL"
$ws2.Range("C9").Value = "BE_IND_8_REAL"
$ws2.Range("E9").Value = "BE_IND_8"
$ws2.Range("F9").Value = "REAL"

$ws2.Range("A10").Value = "CREATE/MODIFY"
$ws2.Range("B10").Value = "BE_IND_9_STRING"
$ws2.Range("C10").Value = "BE_IND_9_STRING"
$ws2.Range("E10").Value = "BE_IND_9"
$ws2.Range("F10").Value = "STRING"

$ws2.Range("A11").Value = "CREATE/MODIFY"
$ws2.Range("B11").Value = "BE_IND_14_REAL"
$ws2.Range("C11").Value = "BE_IND_14_REAL"
$ws2.Range("E11").Value = "BE_IND_14"
$ws2.Range("F11").Value = "REAL"

$ws2.Range("A12").Value = "CREATE/MODIFY"
$ws2.Range("B12").Value = "BE_IND_16_STRING"
$ws2.Range("C12").Value = "BE_IND_16_STRING"
$ws2.Range("E12").Value = "BE_IND_16"
$ws2.Range("F12").Value = "STRING"

$ws2.Range("A13").Value = "CREATE/MODIFY"
$ws2.Range("B13").Value = "BE_IND_34_STRING"
$ws2.Range("C13").Value = "BE_IND_34_STRING"
$ws2.Range("E13").Value = "BE_IND_34"
$ws2.Range("F13").Value = "STRING"

$ws2.Range("A14").Value = "CREATE/MODIFY"
$ws2.Range("B14").Value = "BE_IND_35_REAL"
$ws2.Range("C14").Value = "BE_IND_35_REAL"
$ws2.Range("E14").Value = "BE_IND_35"
$ws2.Range("F14").Value = "REAL"

$ws2.Range("A15").Value = "CREATE/MODIFY"
$ws2.Range("B15").Value = "BE_IND_40_REAL"
$ws2.Range("C15").Value = "BE_IND_40_REAL"
$ws2.Range("E15").Value = "BE_IND_40"
$ws2.Range("F15").Value = "REAL"

$ws2.Range("A16").Value = "CREATE/MODIFY"
$ws2.Range("B16").Value = "BE_IND_44_REAL"
$ws2.Range("C16").Value = "BE_IND_44_REAL"
$ws2.Range("E16").Value = "BE_IND_44"
$ws2.Range("F16").Value = "REAL"

$ws2.Range("A17").Value = "CREATE/MODIFY"
$ws2.Range("B17").Value = "BE_IND_48_REAL"
$ws2.Range("C17").Value = "BE_IND_48_REAL"
$ws2.Range("E17").Value = "BE_IND_48"
$ws2.Range("F17").Value = "REAL"

$ws2.Range("A18").Value = "CREATE/MODIFY"
$ws2.Range("B18").Value = "BE_IND_51_STRING"
$ws2.Range("C18").Value = "BE_IND_51_STRING"
$ws2.Range("E18").Value = "BE_IND_51"
$ws2.Range("F18").Value = "STRING"

$ws2.Range("A19").Value = "CREATE/MODIFY"
$ws2.Range("B19").Value = "BE_IND_55_REAL"
$ws2.Range("C19").Value = "BE_IND_55_REAL"
$ws2.Range("E19").Value = "BE_IND_55"
$ws2.Range("F19").Value = "REAL"

$ws2.Range("A20").Value = "CREATE/MODIFY"
$ws2.Range("B20").Value = "BE_IND_56_REAL"
$ws2.Range("C20").Value = "BE_IND_56_REAL"
$ws2.Range("E20").Value = "BE_IND_56"
$ws2.Range("F20").Value = "REAL"

